$wb = $excel.ActiveWorkbook

# Rename the "Norite" sheet to "Norite Density"
$noriteSheet = $wb.Worksheets.Item("Norite")
$noriteSheet.Name = "Norite Density"

# Make "Norite Density" the active/selected sheet (instead of "Attenuation Coefficients")
$noriteSheet.Activate()

$wb.Save()
